# CurrentRelay workbook: replace the boolean "IsTopBuyed"/"IsNew" flag
# columns with a single new "DisplayClass" column (commit: "Added Display
# class").
#
# Column G currently holds "IsTopBuyed" and column H holds "IsNew" (both
# empty for the existing data row). We delete both columns - which shifts
# every later column two places to the left - then insert one fresh blank
# column back in at G and label it "DisplayClass". Net effect: one column
# fewer overall, and everything from (old) "Type" onward now sits one
# column to the left of where it used to be.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CurrentRelay")

# Remove the two retired flag columns (G = IsTopBuyed, H = IsNew).
# Deleting column 7 twice removes both, because the second delete operates
# on whatever has shifted into column 7 after the first delete.
$ws1.Columns.Item(7).Delete()
$ws1.Columns.Item(7).Delete()

# Insert a fresh empty column at G for the new field and label it.
$ws1.Columns.Item(7).Insert()
$ws1.Cells.Item(1, 7).Value = "DisplayClass"

# Match the author's recorded selection on save.
$ws1.Range("G2").Select()
